# Auto-generated Excel COM-interop script to apply cell value edits
# per the unified diff of Jenova_Profits.xlsx canonical OOXML.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2447.36
$ws.Range("I38").Value = 2447.36
$ws.Range("K38").Value = 7342.08
$ws.Range("M38").Value = -6970.08
$ws.Range("H51").Value = 3464.8333
$ws.Range("I51").Value = 3322.25
$ws.Range("K51").Value = 3322.25
$ws.Range("M51").Value = -2838.25
$ws.Range("H58").Value = 3575
$ws.Range("J58").Value = 7307.6924
$ws.Range("L58").Value = 21923.0772
$ws.Range("N58").Value = -22223.0772
$ws.Range("H80").Value = 4228.643
$ws.Range("I80").Value = 6683.4116
$ws.Range("J80").Value = 434.9091
$ws.Range("K80").Value = 20050.2348
$ws.Range("L80").Value = 1304.7273
$ws.Range("M80").Value = -19052.2348
$ws.Range("N80").Value = -3300.7273
$ws.Range("H83").Value = 4228.643
$ws.Range("I83").Value = 6683.4116
$ws.Range("J83").Value = 434.9091
$ws.Range("K83").Value = 60150.7044
$ws.Range("L83").Value = 3914.1819
$ws.Range("M83").Value = -55158.7044
$ws.Range("N83").Value = -13898.1819
$ws.Range("H88").Value = 1103.1666
$ws.Range("J88").Value = 1697
$ws.Range("L88").Value = 1697
$ws.Range("N88").Value = -2509
$ws.Range("H91").Value = 1103.1666
$ws.Range("J91").Value = 1697
$ws.Range("L91").Value = 1697
$ws.Range("N91").Value = -4505
$ws.Range("H112").Value = 1983.8462
$ws.Range("J112").Value = 1999.1666
$ws.Range("L112").Value = 5997.4998
$ws.Range("N112").Value = -8213.4998
$ws.Range("H137").Value = 9186.727999999999
$ws.Range("I137").Value = 8309.6
$ws.Range("J137").Value = 9917.666999999999
$ws.Range("K137").Value = 24928.8
$ws.Range("L137").Value = 29753.001
$ws.Range("M137").Value = -22378.8
$ws.Range("N137").Value = -34853.001
$ws.Range("H138").Value = 9127.841
$ws.Range("J138").Value = 9834.233
$ws.Range("L138").Value = 29502.699
$ws.Range("N138").Value = -39782.699
$ws.Range("H141").Value = 6166.7334
$ws.Range("I141").Value = 6166.7334
$ws.Range("K141").Value = 18500.2002
$ws.Range("M141").Value = -13320.2002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 75208.53
$ws.Range("I2").Value = 138851.25
$ws.Range("J2").Value = 2474
$ws.Range("K2").Value = 138851.25
$ws.Range("L2").Value = 2474
$ws.Range("M2").Value = -138738.25
$ws.Range("N2").Value = -2700
$ws.Range("H32").Value = 22602.508
$ws.Range("I32").Value = 13119.97
$ws.Range("J32").Value = 33033.3
$ws.Range("K32").Value = 13119.97
$ws.Range("L32").Value = 33033.3
$ws.Range("M32").Value = -12832.97
$ws.Range("N32").Value = -33607.3
$ws.Range("H74").Value = 4823.4546
$ws.Range("I74").Value = 3882.375
$ws.Range("K74").Value = 3882.375
$ws.Range("M74").Value = -3008.375
$ws.Range("H77").Value = 4823.4546
$ws.Range("I77").Value = 3882.375
$ws.Range("K77").Value = 19411.875
$ws.Range("M77").Value = -15043.875
$ws.Range("H102").Value = 2251.516
$ws.Range("I102").Value = 2253.3215
$ws.Range("K102").Value = 2253.3215
$ws.Range("M102").Value = -631.3215
$ws.Range("H116").Value = 75208.53
$ws.Range("I116").Value = 138851.25
$ws.Range("J116").Value = 2474
$ws.Range("K116").Value = 138851.25
$ws.Range("L116").Value = 2474
$ws.Range("M116").Value = -136557.25
$ws.Range("N116").Value = -7062
$ws.Range("H122").Value = 4833.3335
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4833.3335
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 14500.0005
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -19400.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 75208.53
$ws.Range("I3").Value = 138851.25
$ws.Range("J3").Value = 2474
$ws.Range("K3").Value = 138851.25
$ws.Range("L3").Value = 2474
$ws.Range("M3").Value = -138737.25
$ws.Range("N3").Value = -2702
$ws.Range("H59").Value = 99990
$ws.Range("J59").Value = 99990
$ws.Range("L59").Value = 99990
$ws.Range("N59").Value = -101684
$ws.Range("H94").Value = 2530.5
$ws.Range("I94").Value = 2077.5293
$ws.Range("J94").Value = 3230.5454
$ws.Range("K94").Value = 2077.5293
$ws.Range("L94").Value = 3230.5454
$ws.Range("M94").Value = -1626.5293
$ws.Range("N94").Value = -4132.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 74191
$ws.Range("I31").Value = 2974.923
$ws.Range("J31").Value = 1000000
$ws.Range("K31").Value = 2974.923
$ws.Range("L31").Value = 1000000
$ws.Range("M31").Value = -2679.923
$ws.Range("N31").Value = -1000590
$ws.Range("H34").Value = 74191
$ws.Range("I34").Value = 2974.923
$ws.Range("J34").Value = 1000000
$ws.Range("K34").Value = 2974.923
$ws.Range("L34").Value = 1000000
$ws.Range("M34").Value = -2772.923
$ws.Range("N34").Value = -1000404
$ws.Range("H58").Value = 2485.7144
$ws.Range("I58").Value = 2262.875
$ws.Range("J58").Value = 2782.8333
$ws.Range("K58").Value = 2262.875
$ws.Range("L58").Value = 2782.8333
$ws.Range("M58").Value = -2059.875
$ws.Range("N58").Value = -3188.8333
$ws.Range("H82").Value = 75248
$ws.Range("J82").Value = 75248
$ws.Range("L82").Value = 75248
$ws.Range("N82").Value = -75970
$ws.Range("H85").Value = 75248
$ws.Range("J85").Value = 75248
$ws.Range("L85").Value = 75248
$ws.Range("N85").Value = -77744
$ws.Range("H105").Value = 807.375
$ws.Range("I105").Value = 807.375
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 807.375
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 939.625
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 530248
$ws.Range("I134").Value = 4652.0713
$ws.Range("J134").Value = 2001916.6
$ws.Range("K134").Value = 13956.2139
$ws.Range("L134").Value = 6005749.800000001
$ws.Range("M134").Value = -11421.2139
$ws.Range("N134").Value = -6010819.800000001
$ws.Range("H136").Value = 2485.7144
$ws.Range("I136").Value = 2262.875
$ws.Range("J136").Value = 2782.8333
$ws.Range("K136").Value = 6788.625
$ws.Range("L136").Value = 8348.499899999999
$ws.Range("M136").Value = -4238.625
$ws.Range("N136").Value = -13448.4999
$ws.Range("H141").Value = 468231
$ws.Range("J141").Value = 510054.7
$ws.Range("L141").Value = 510054.7
$ws.Range("N141").Value = -520414.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 977802.3
$ws.Range("I4").Value = 1052560.6
$ws.Range("J4").Value = 379736
$ws.Range("K4").Value = 3157681.8
$ws.Range("L4").Value = 1139208
$ws.Range("M4").Value = -3157569.8
$ws.Range("N4").Value = -1139432
$ws.Range("H34").Value = 129062.625
$ws.Range("J34").Value = 206340
$ws.Range("L34").Value = 619020
$ws.Range("N34").Value = -619188
$ws.Range("H39").Value = 13887.35
$ws.Range("I39").Value = 5178.4287
$ws.Range("J39").Value = 18576.77
$ws.Range("K39").Value = 15535.2861
$ws.Range("L39").Value = 55730.31
$ws.Range("M39").Value = -15241.2861
$ws.Range("N39").Value = -56318.31
$ws.Range("H55").Value = 9428.571
$ws.Range("J55").Value = 10916.667
$ws.Range("L55").Value = 32750.001
$ws.Range("N55").Value = -33104.001
$ws.Range("H140").Value = 4859.5264
$ws.Range("I140").Value = 4069.75
$ws.Range("K140").Value = 12209.25
$ws.Range("M140").Value = -7029.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 38692.137
$ws.Range("J24").Value = 38692.137
$ws.Range("L24").Value = 38692.137
$ws.Range("N24").Value = -39038.137
$ws.Range("H25").Value = 2806.3333
$ws.Range("J25").Value = 2806.3333
$ws.Range("L25").Value = 2806.3333
$ws.Range("N25").Value = -3864.3333
$ws.Range("I122").Value = 3678
$ws.Range("J122").Value = 4233.3335
$ws.Range("K122").Value = 11034
$ws.Range("L122").Value = 12700.0005
$ws.Range("M122").Value = -8584
$ws.Range("N122").Value = -17600.0005
$ws.Range("H123").Value = 80498.5
$ws.Range("J123").Value = 80498.5
$ws.Range("L123").Value = 80498.5
$ws.Range("N123").Value = -85398.5
$ws.Range("H132").Value = 70674.94500000001
$ws.Range("I132").Value = 11566.363
$ws.Range("J132").Value = 163559.86
$ws.Range("K132").Value = 34699.089
$ws.Range("L132").Value = 490679.58
$ws.Range("M132").Value = -32169.089
$ws.Range("N132").Value = -495739.58

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1975
$ws.Range("I22").Value = 1975
$ws.Range("K22").Value = 1975
$ws.Range("M22").Value = -1680
$ws.Range("H27").Value = 1975
$ws.Range("I27").Value = 1975
$ws.Range("K27").Value = 1975
$ws.Range("M27").Value = -1868
$ws.Range("H46").Value = 3366.889
$ws.Range("I46").Value = 3125
$ws.Range("J46").Value = 3560.4
$ws.Range("K46").Value = 3125
$ws.Range("L46").Value = 3560.4
$ws.Range("M46").Value = -2937
$ws.Range("N46").Value = -3936.4
$ws.Range("H76").Value = 14183
$ws.Range("J76").Value = 16144
$ws.Range("L76").Value = 16144
$ws.Range("N76").Value = -16820
$ws.Range("H79").Value = 14183
$ws.Range("J79").Value = 16144
$ws.Range("L79").Value = 16144
$ws.Range("N79").Value = -18484
$ws.Range("H132").Value = 7025.4546
$ws.Range("I132").Value = 6896.5713
$ws.Range("K132").Value = 20689.7139
$ws.Range("M132").Value = -18159.7139
$ws.Range("H133").Value = 56444.223
$ws.Range("J133").Value = 56444.223
$ws.Range("L133").Value = 56444.223
$ws.Range("N133").Value = -61504.223

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 10000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 10000
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -10224
$ws.Range("H132").Value = 103965.25
$ws.Range("I132").Value = 4286.4443
$ws.Range("K132").Value = 12859.3329
$ws.Range("M132").Value = -10329.3329
